# Update the "想去人数" (want-to-go count) figures in column F for the
# rows that changed between the previous build and the one generated at
# commit 456a3b4. The same underlying data is duplicated on the "展览"
# sheet and the "全部类型" sheet, so both need to be updated identically.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1832
    $ws.Range("F3").Value = 8365
    $ws.Range("F5").Value = 353
}
